$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.272.46"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +2.24%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'1.877.94"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +1.48%  "
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'  +0.16%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'315.85"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +0.63%  "
$ws.Range("E5").ClearFormats()
$ws.Range("E6").Value = "'  +0.23%  "
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'0.4309"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  +1.57%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.3731"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  +2.44%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.07401"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  +1.25%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.8851"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  +0.80%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'21.14"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  +1.72%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'1.966.30"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  +6.94%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'5.489"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +2.70%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'6.634"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +1.46%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'0.06982"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  +1.14%  "
$ws.Range("E15").ClearFormats()
$ws.Range("E16").Value = "'  +0.26%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'81.29"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +2.44%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'0.000009118"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +2.36%  "
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = "'  +0.31%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'15.64"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  +1.59%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'28.348.81"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  +2.50%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'5.099"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +2.22%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'10.94"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  +4.71%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'2.229.36"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  +7.61%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'1.969"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  -0.65%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'154.51"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  +1.15%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'18.76"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  -1.09%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'5.384"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  +2.04%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'117.19"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  -3.46%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'1.862"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  -1.32%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'0.08986"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  +1.09%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'0.7900"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  +2.83%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'4.678"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  +2.28%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'1.183"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  +6.99%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'2.957"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  -0.51%  "
$ws.Range("E35").ClearFormats()
$ws.Range("E36").Value = "'  +0.30%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'1.128"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  +2.86%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'0.05449"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  +1.51%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'0.01965"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  +1.55%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'2.892"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  +2.71%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.5165"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  +0.99%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'0.1683"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  +1.84%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'6.881"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  -0.25%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'8.648"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  +4.00%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'10.57"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  +1.95%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'0.06603"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  +0.95%  "
$ws.Range("E46").ClearFormats()
$ws.Range("E47").Value = "'  +0.13%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'106.31"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +1.34%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'1.002"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +0.29%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'1.654"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  +1.61%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'1.822"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  +5.11%  "
$ws.Range("E51").ClearFormats()
